$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Duplicate the existing "2022-Q2" sheet (placing the copy right after it) so
# the old quarter's data is preserved on its own tab, then turn the original
# tab into the new "2022-Q3" quarter.
$ws2.Copy($null, $ws2)
$wsCopy = $wb.Worksheets.Item(3)
$ws2.Name = "2022-Q3"
$wsCopy.Name = "2022-Q2"

# Reuse the "总计" sheet's header / leading-column formatting for the new
# quarter's table before the old figures are overwritten.
$ws1.Range("B1:D1").Copy()
$ws2.Range("B1:H1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws2.Range("A2:A4").PasteSpecial(-4122)

# Replace the fund holdings with the 2022-Q3 figures.
$ws2.Range("A2:H4").ClearContents()

$ws2.Range("A2").Value = 0
$ws2.Range("C2").Value = "兴业安保优选混合"
$ws2.Range("H2").Value = 8

$ws2.Range("A3").Value = 1
$ws2.Range("C3").Value = "景顺长城改革机遇灵活配置混合A"
$ws2.Range("H3").Value = 8

$ws2.Range("A4").Value = 2
$ws2.Range("C4").Value = "景顺长城改革机遇灵活配置混合C"
$ws2.Range("H4").Value = 8

# B and D:G hold numeric-looking figures that must stay text, not numbers -
# force a text format before assigning so Excel doesn't auto-convert them,
# then drop back to the Normal style so no stray number format lingers.
$ws2.Range("B2:B4").NumberFormat = "@"
$ws2.Range("D2:G4").NumberFormat = "@"

$ws2.Range("B2").Value = "006366"
$ws2.Range("D2").Value = "1.26"
$ws2.Range("E2").Value = "87.20"
$ws2.Range("F2").Value = "4.77"
$ws2.Range("G2").Value = "0.0601"

$ws2.Range("B3").Value = "001535"
$ws2.Range("D3").Value = "0.26"
$ws2.Range("E3").Value = "64.66"
$ws2.Range("F3").Value = "2.42"
$ws2.Range("G3").Value = "0.0063"

$ws2.Range("B4").Value = "007945"
$ws2.Range("D4").Value = "0.06"
$ws2.Range("E4").Value = "64.66"
$ws2.Range("F4").Value = "2.42"
$ws2.Range("G4").Value = "0.0015"

$ws2.Range("B2:B4").Style = "Normal"
$ws2.Range("D2:G4").Style = "Normal"

# Update the "总计" (summary) sheet: push the existing 2022-Q2 summary row
# down to row 3, then overwrite row 2 with the new 2022-Q3 summary.
$ws1.Range("A2:D2").Copy($ws1.Range("A3:D3"))
$ws1.Range("A3").Value = 1

$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 3
$ws1.Range("D2").Value = 0.07000000000000001

# Keep the original active-sheet selection (the "总计" overview tab).
$ws1.Activate()
